# Organizational_Roles_Responsibilities.docx
# Commit: "Updated roles and responsibilities for System Engineering
#          Removed Bill Reed, as he is working on GOLF."
#
# Net change:
#   1. Drop "and Bill Reed" from the System engineering roles/responsibilities
#      sentence so it reads "...(Jonathan Black for Spacecraft TBD for antennas)".
#   2. Word maintains a single hidden "_GoBack" bookmark marking the most
#      recent edit location. Because this edit lands earlier in the document
#      than where "_GoBack" previously sat (end of doc), the bookmark moves
#      to the empty paragraph that immediately follows the edit.

$d = $word.ActiveDocument

# --- Relocate the "_GoBack" bookmark ------------------------------------
# Remove it from its old spot (end of document, after the "antennas))" text).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Re-create it at the empty paragraph right after "...Tom Clark, Marc Franco)",
# which is where Word leaves it after this edit.
$anchor = $d.Content
$anchor.Find.Execute(
    "Tom Clark, Marc Franco)", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null
$targetPara = $anchor.Paragraphs(1).Next()
$d.Bookmarks.Add("_GoBack", $targetPara.Range) | Out-Null

# --- Text edit: remove "and Bill Reed" ----------------------------------
$rng = $d.Content
$rng.Find.Execute(
    "and Bill Reed ", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 2) | Out-Null

Write-Output "Bill Reed removed; _GoBack relocated."
